$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DQ_Metrics")

# Insert a new column before column J (shifts J:T -> K:U)
$ws.Range("J1").EntireColumn.Insert()

# --- Header row (row 1) updates ---
$ws.Range("C1").Value = "item_completeness_rate"
$ws.Range("D1").Value = "value_completeness_rate"
$ws.Range("F1").Value = "range_plausibility_rate"
$ws.Range("I1").Value = "rdCase_dissimilarity_rate"
$ws.Range("J1").Value = "rdCase_rel_py_ipat"

# New trailing header columns V1:AB1
$ws.Range("V1").Value = "missing_item_no_py"
$ws.Range("W1").Value = "missing_value_no_py"
$ws.Range("X1").Value = "orphaMissing_no_py"
$ws.Range("Y1").Value = "implausible_codeLink_no_py"
$ws.Range("Z1").Value = "outlier_no_py"
$ws.Range("AA1").Value = "ambigous_rdCase_no_py"
$ws.Range("AB1").Value = "duplicateRdCase_no_py"

# Apply header style (bold/centered, matches style index 1) to new cells
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108

$ws.Range("V1:AB1").Font.Bold = $true
$ws.Range("V1:AB1").HorizontalAlignment = -4108

# --- Data row (row 2) updates ---
$ws.Range("C2").Value = 78.56999999999999
$ws.Range("D2").Value = 96.22
$ws.Range("F2").Value = 99.83
$ws.Range("I2").Value = 97
$ws.Range("J2").Value = 1

# New trailing data columns V2:AB2
$ws.Range("V2").Value = 3
$ws.Range("W2").Value = 518
$ws.Range("X2").Value = 11
$ws.Range("Y2").Value = 22
$ws.Range("Z2").Value = 8
$ws.Range("AA2").Value = 25
$ws.Range("AB2").Value = 3
